$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 148, shifting existing rows 148:182 down to 149:183
$ws.Rows("148:148").Insert()

# Populate the newly inserted row 148 with the new weekly price-report entry
$ws.Range("A148").Value = 10
$ws.Range("B148").Value = "Vega Modelo de Temuco"
$ws.Range("C148").Value = "La Araucanía"
$ws.Range("D148").Value = 44782
$ws.Range("E148").Value = 9
$ws.Range("F148").Value = 100114007
$ws.Range("G148").Value = "Jengibre"
$ws.Range("H148").Value = "Sin especificar"
$ws.Range("I148").Value = "Primera"
$ws.Range("J148").Value = 140
$ws.Range("K148").Value = 16000
$ws.Range("L148").Value = 20000
$ws.Range("M148").Value = 17429
$ws.Range("N148").Value = "$/caja 13 kilos"
$ws.Range("O148").Value = "Perú"
$ws.Range("P148").Value = 1341
$ws.Range("Q148").Value = 13
$ws.Range("R148").Value = "Hortaliza"
